# "Natmi following Dr Hou advice"
# The sending/target cluster set grows from {ECs, FAPs} to {ECs, FAPs, sCs},
# so the 2x2 (4-row) ligand/receptor table becomes a 3x3 (9-row) table with
# recomputed specificity metrics. Rewrite A2:T10 in one shot with the full
# new table (row 1 headers are untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20

$data[0,0] = "ECs"
$data[0,1] = "Il23a"
$data[0,2] = "Il12rb1"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 31.66883366666667
$data[0,7] = 95.006501
$data[0,8] = 0.9615169947075755
$data[0,9] = 0.9615169947075756
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.1660766666666667
$data[0,13] = 0.49823
$data[0,14] = 0.1038853751668853
$data[0,15] = 0.1038853751668853
$data[0,16] = 5.259454332581112
$data[0,17] = 47.33508899323
$data[0,18] = 0.09988755372453255
$data[0,19] = 0.09988755372453255

$data[1,0] = "ECs"
$data[1,1] = "Il23a"
$data[1,2] = "Il12rb1"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 31.66883366666667
$data[1,7] = 95.006501
$data[1,8] = 0.9615169947075755
$data[1,9] = 0.9615169947075756
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 1.211564666666667
$data[1,13] = 3.634694
$data[1,14] = 0.7578659450591634
$data[1,15] = 0.7578659450591634
$data[1,16] = 38.36883990507711
$data[1,17] = 345.319559145694
$data[1,18] = 0.7287009858845033
$data[1,19] = 0.7287009858845034

$data[2,0] = "ECs"
$data[2,1] = "Il23a"
$data[2,2] = "Il12rb1"
$data[2,3] = "sCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 31.66883366666667
$data[2,7] = 95.006501
$data[2,8] = 0.9615169947075755
$data[2,9] = 0.9615169947075756
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.2210116666666666
$data[2,13] = 0.6630349999999999
$data[2,14] = 0.1382486797739514
$data[2,15] = 0.1382486797739514
$data[2,16] = 6.999181710059443
$data[2,17] = 62.99263539053499
$data[2,18] = 0.1329284550985397
$data[2,19] = 0.1329284550985397

$data[3,0] = "FAPs"
$data[3,1] = "Il23a"
$data[3,2] = "Il12rb1"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 1.032613
$data[3,7] = 3.097839
$data[3,8] = 0.03135180028751844
$data[3,9] = 0.03135180028751844
$data[3,10] = 2
$data[3,11] = 0.6666666666666666
$data[3,12] = 0.1660766666666667
$data[3,13] = 0.49823
$data[3,14] = 0.1038853751668853
$data[3,15] = 0.1038853751668853
$data[3,16] = 0.1714929249966667
$data[3,17] = 1.54343632497
$data[3,18] = 0.003256993535026116
$data[3,19] = 0.003256993535026116

$data[4,0] = "FAPs"
$data[4,1] = "Il23a"
$data[4,2] = "Il12rb1"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1.032613
$data[4,7] = 3.097839
$data[4,8] = 0.03135180028751844
$data[4,9] = 0.03135180028751844
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.211564666666667
$data[4,13] = 3.634694
$data[4,14] = 0.7578659450591634
$data[4,15] = 0.7578659450591634
$data[4,16] = 1.251077425140667
$data[4,17] = 11.259696826266
$data[4,18] = 0.02376046175420631
$data[4,19] = 0.02376046175420631

$data[5,0] = "FAPs"
$data[5,1] = "Il23a"
$data[5,2] = "Il12rb1"
$data[5,3] = "sCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.032613
$data[5,7] = 3.097839
$data[5,8] = 0.03135180028751844
$data[5,9] = 0.03135180028751844
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.2210116666666666
$data[5,13] = 0.6630349999999999
$data[5,14] = 0.1382486797739514
$data[5,15] = 0.1382486797739514
$data[5,16] = 0.2282195201516666
$data[5,17] = 2.053975681365
$data[5,18] = 0.004334344998286012
$data[5,19] = 0.004334344998286012

$data[6,0] = "sCs"
$data[6,1] = "Il23a"
$data[6,2] = "Il12rb1"
$data[6,3] = "ECs"
$data[6,4] = 2
$data[6,5] = 0.6666666666666666
$data[6,6] = 0.2348756666666667
$data[6,7] = 0.704627
$data[6,8] = 0.007131205004906082
$data[6,9] = 0.007131205004906083
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.1660766666666667
$data[6,13] = 0.49823
$data[6,14] = 0.1038853751668853
$data[6,15] = 0.1038853751668853
$data[6,16] = 0.03900736780111112
$data[6,17] = 0.35106631021
$data[6,18] = 0.0007408279073266385
$data[6,19] = 0.0007408279073266386

$data[7,0] = "sCs"
$data[7,1] = "Il23a"
$data[7,2] = "Il12rb1"
$data[7,3] = "FAPs"
$data[7,4] = 2
$data[7,5] = 0.6666666666666666
$data[7,6] = 0.2348756666666667
$data[7,7] = 0.704627
$data[7,8] = 0.007131205004906082
$data[7,9] = 0.007131205004906083
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 1.211564666666667
$data[7,13] = 3.634694
$data[7,14] = 0.7578659450591634
$data[7,15] = 0.7578659450591634
$data[7,16] = 0.2845670587931111
$data[7,17] = 2.561103529138
$data[7,18] = 0.005404497420453784
$data[7,19] = 0.005404497420453785

$data[8,0] = "sCs"
$data[8,1] = "Il23a"
$data[8,2] = "Il12rb1"
$data[8,3] = "sCs"
$data[8,4] = 2
$data[8,5] = 0.6666666666666666
$data[8,6] = 0.2348756666666667
$data[8,7] = 0.704627
$data[8,8] = 0.007131205004906082
$data[8,9] = 0.007131205004906083
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.2210116666666666
$data[8,13] = 0.6630349999999999
$data[8,14] = 0.1382486797739514
$data[8,15] = 0.1382486797739514
$data[8,16] = 0.05191026254944444
$data[8,17] = 0.4671923629449999
$data[8,18] = 0.0009858796771256603
$data[8,19] = 0.0009858796771256603

$ws.Range("A2:T10").Value = $data
